$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns I (I0) and J (IF), matching the style of the
# existing header row (style used by B1:H1) by copying H1's formatting over.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new columns for each data row (rows 2-36).
# Column J mirrors column H's value for every row; column I is 1 for every
# row except row 33, where I=6 and J=9 (and H stays unchanged).
for ($r = 2; $r -le 36; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2

    if ($r -eq 33) {
        $ws.Cells.Item($r, 9).Value = 6
        $ws.Cells.Item($r, 10).Value = 9
    } else {
        $ws.Cells.Item($r, 9).Value = 1
        $ws.Cells.Item($r, 10).Value = $hVal
    }
}
